$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing last column (Q) into the new column (R)
# for every data row, then fill in the 2021 values.
$ws.Range("Q4:Q34").Copy() | Out-Null
$ws.Range("R4").PasteSpecial(-4122) | Out-Null

$ws.Range("R4").Value = 2021
$ws.Range("R5").Value = 11.9
$ws.Range("R6").Value = 13.1
$ws.Range("R7").Value = 10.6
$ws.Range("R8").Value = 11
$ws.Range("R9").Value = 10
$ws.Range("R10").Value = 12
$ws.Range("R11").Value = 10.199999999999999
$ws.Range("R12").Value = 10.5
$ws.Range("R13").Value = 10
$ws.Range("R14").Value = 19.399999999999999
$ws.Range("R15").Value = 22.3
$ws.Range("R16").Value = 16.399999999999999
$ws.Range("R17").Value = 9.4
$ws.Range("R18").Value = 11.4
$ws.Range("R19").Value = 7.3
$ws.Range("R20").Value = 3.1
$ws.Range("R21").Value = 2.9
$ws.Range("R22").Value = 3.4
$ws.Range("R23").Value = 15
$ws.Range("R24").Value = 17.3
$ws.Range("R25").Value = 12.7
$ws.Range("R26").Value = 7.9
$ws.Range("R27").Value = 8.4
$ws.Range("R28").Value = 7.4
$ws.Range("R29").Value = 15.2
$ws.Range("R30").Value = 17.600000000000001
$ws.Range("R31").Value = 12.6
$ws.Range("R32").Value = 27.9
$ws.Range("R33").Value = 32.700000000000003
$ws.Range("R34").Value = 22.8

# Match the author's new selection (cell R3, just above the new data)
$ws.Range("R3").Select() | Out-Null
